# Variable packaging consignment class, sample data, and argument for config
#
# - Sheet2 gains a new "variable packaging" File type option (Table2, B4),
#   which is now covered by the B14 dropdown list on Sheet1.
# - Sheet1's "consignment generation method" (B6) switches from the
#   parameter-based default to the input_file method.
# - Sheet1's "File name" (B15) points at the new sample CSV for the
#   variable-packaging example instead of the AQIM sample.
# - Sheet1's "sample size method" (B38) switches from proportion to
#   hypergeometric.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Add the new "variable packaging" file-type option to Sheet2's lookup table.
$ws2.Range("B4").Value = "variable packaging"
$ws2.Columns.Item(2).AutoFit()

# Update the example parameter values on Sheet1.
$ws1.Range("B6").Value = "input_file"
$ws1.Range("B15").Value = '"data/varipack_sample.csv"'
$ws1.Range("B38").Value = "hypergeometric"

# Leave the same selection state the saved workbook ends up with.
$ws2.Activate() | Out-Null
$ws2.Range("C12").Select() | Out-Null
$ws1.Activate() | Out-Null
$ws1.Range("C10").Select() | Out-Null
